# Reorganize the TODO list into category (col B) / task (col C) rows,
# drop the finished "Licencja do Pycharma" item, add the new Arduino task,
# and flip the detail-column alignment from right to left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Break the old B5:B6 merge before rewriting the grid --------------------
$ws.Range("B5:B6").UnMerge()

# --- Category column (B) -----------------------------------------------------
$ws.Range("B2").Value = "Audio"
$ws.Range("B3").Value = ""
$ws.Range("B4").Value = "Silniki krokowe"
$ws.Range("B5").Value = ""
$ws.Range("B6").Value = ""

# --- Detail column (C) --------------------------------------------------------
$ws.Range("C2").Value = "Uszkodzony układ audio u Piotera"
$ws.Range("C3").Value = "Użyć nowych odzyskanych głośników"
$ws.Range("C4").Value = "Skonfigurować nową płytkę arduino+drv8825 zamiast ender3"
$ws.Range("C5").Value = "Przełożenie zębatek wymaga aktualizacji"
$ws.Range("C6").Value = "Poeksperymentować z trybami pracy driverów w silnikach"

# --- New merges: one category cell per task group -----------------------------
$ws.Range("B2:B3").Merge()
$ws.Range("B4:B6").Merge()

# --- Row height (row 4 now wraps a long task, same as rows 5/6) ---------------
$ws.Rows.Item(4).RowHeight = 30

# --- Column widths (slightly widened) ------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 35.67
$ws.Columns.Item(3).ColumnWidth = 33.67

# --- Detail column for the Audio rows switches from right- to left-aligned ----
$r = $ws.Range("C2:C3")
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.WrapText = $false

# --- Selection moves up one row now that the sheet is one row shorter ---------
$null = $ws.Range("C7").Select()
